# Add a "Save" column (H) to the sheet, mirroring the header style used by
# the other header cells (B1:G1), then populate H2:H39 with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from G1 (bold/bordered/centered header style) onto H1,
# then overwrite the copied value with the new header text "Save".
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New "Save" column values for rows 2-39 (one entry per row, in order).
$saveValues = @(0,1,0,0,1,0,1,0,0,0,0,1,1,0,1,0,1,0,0,0,0,0,0,1,0,0,1,1,0,0,1,0,0,1,1,1,1,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
